# Add "STEPS" tasks/notes columns (E:H) and additional task rows to the
# Tasks sheet, matching the author's upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "STEPS" header + notes columns -----------------------------------
$ws.Range("E1").Value = "STEPS"

$ws.Range("E3").Value = "install inkspace"
$ws.Range("F3").Value = "try first example"
$ws.Range("G3").Value = "apply to my bus"
$ws.Range("H3").Value = "transform to schemas"

$ws.Range("E14:F14").Merge()
$ws.Range("E14").Value = "clean data using marcos"
$ws.Range("G14").Value = "create macros"
$ws.Range("H14").Value = "transform to addons"

# --- Rename row 17's task, then insert the original "car wash" task back
#     as a new row 18 (pushes the old red/bold "ubuntu shutdown" row down
#     to row 19) -------------------------------------------------------
$ws.Range("B17").Value = "go beach"

$ws.Rows.Item(18).Insert()
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "car wash"
$ws.Range("C18").Value = 1

# --- Row 19 used to be the old "ubuntu shutdown" row (18, 0) with a
#     special red/bold style; give it the new content and plain styling
#     used by the rest of the table ------------------------------------
$ws.Range("B19").Value = "git"
$ws.Range("C19").Value = 1
$ws.Range("E19").Value = "learn about"
$ws.Range("F19").Value = "tasks example"

$ws.Range("A19:F19").Font.Name = "Calibri"
$ws.Range("A19:F19").Font.Size = 11
$ws.Range("A19:F19").Font.Bold = $false
$ws.Range("A19:F19").Font.Color = 0
$ws.Rows.Item(19).RowHeight = 15

# --- New trailing rows ------------------------------------------------
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "travaux jardin"
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = "Vendredi"

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "travaux maison"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = "Vendredi"

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Ubuntu 16.04"
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = "Daily"
$ws.Range("E22").Value = "until get solved"

$ws.Range("A22").HorizontalAlignment = -4108
$ws.Rows.Item(22).RowHeight = 13.8

# --- Column widths for the new columns, matching the author's manual
#     resize -------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 16.8333333333333
$ws.Columns.Item(5).ColumnWidth = 14
$ws.Columns.Item(6).ColumnWidth = 14.8333333333333
$ws.Columns.Item(7).ColumnWidth = 14.3333333333333
$ws.Columns.Item(8).ColumnWidth = 19.6666666666667

$ws.Range("E22").Select()
